# Loan RBI, Variable Instalments
#
# The "Repayment Schedule" sheet gains a new (empty) column N, pushing the
# existing "Late" (old column N) and "Outstanding" (old column P) columns
# one slot to the right (to O and Q respectively). The active/selected
# sheet also moves from "NewLoanInput" to "Repayment Schedule", with a new
# selection on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before the old "Late" column (column N / 14),
# shifting "Late" -> O and "Outstanding" -> Q.
$ws.Columns.Item(14).Insert()

# The newly inserted column picks up the width of the column just to its
# left (column M), matching Excel's normal "insert column" behaviour.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth()

# Update the selection on the Repayment Schedule sheet...
$ws.Range("J18").Select() | Out-Null

# ...and make it the active (selected) tab of the workbook, which moves
# tabSelected off "NewLoanInput" and onto "Repayment Schedule".
$ws.Activate() | Out-Null
